# Update ticket/price-related counter figures (column F) on the
# "展览" and "全部类型" worksheets to reflect the latest scraped totals.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 1708
$ws1.Range("F11").Value = 1764
$ws1.Range("F13").Value = 107
$ws1.Range("F14").Value = 419
$ws1.Range("F16").Value = 201
$ws1.Range("F21").Value = 656

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 1708
$ws4.Range("F12").Value = 1764
$ws4.Range("F14").Value = 107
$ws4.Range("F15").Value = 419
$ws4.Range("F17").Value = 201
$ws4.Range("F22").Value = 656
